$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename sheet
$ws.Name = "P1 - Historical"

# 2. Copy the header style (bold, border, centered) from H1 into the new I1:T1 header cells
$ws.Range("H1").Copy()
$ws.Range("I1:T1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 3. Update existing header text (drop "- No costs" suffix) and add new headers
$ws.Range("C1").Value = "EW"
$ws.Range("E1").Value = "MVP"
$ws.Range("G1").Value = "MSR"
$ws.Range("I1").Value = "Min. risk CVaR"
$ws.Range("J1").Value = "Min. risk CVaR - With costs"
$ws.Range("K1").Value = "Optimal CVaR"
$ws.Range("L1").Value = "Optimal CVaR - With costs"
$ws.Range("M1").Value = "Min. risk CDaR"
$ws.Range("N1").Value = "Min. risk CDaR - With costs"
$ws.Range("O1").Value = "Optimal CDaR"
$ws.Range("P1").Value = "Optimal CDaR - With costs"
$ws.Range("Q1").Value = "Min. risk Omega"
$ws.Range("R1").Value = "Min. risk Omega - With costs"
$ws.Range("S1").Value = "Optimal Omega"
$ws.Range("T1").Value = "Optimal Omega - With costs"

# 4. Round the FoF (benchmark) column values to 6 decimal places
$ws.Range("B2").Value = 0.049341
$ws.Range("B3").Value = 0.049341
$ws.Range("B4").Value = 0.056391
$ws.Range("B5").Value = 0.222035
$ws.Range("B6").Value = 0.072119
$ws.Range("B7").Value = 0.242259
$ws.Range("B8").Value = 0.122638
$ws.Range("B9").Value = 0.031147
$ws.Range("B11").Value = 0.647835
$ws.Range("B12").Value = 0.360859

# 5. Populate the new strategy columns (I:T) for rows 2-14
$ws.Range("I2").Value = 0.06695885754874564
$ws.Range("J2").Value = 0.06321119010299214
$ws.Range("K2").Value = 0.07798306381709619
$ws.Range("L2").Value = 0.03916420041440416
$ws.Range("M2").Value = 0.06474864412991455
$ws.Range("N2").Value = 0.06073028485893994
$ws.Range("O2").Value = 0.07248667530731501
$ws.Range("P2").Value = 0.03690580701796131
$ws.Range("Q2").Value = 0.06635623267305446
$ws.Range("R2").Value = 0.05645993002253635
$ws.Range("S2").Value = 0.06846526081957172
$ws.Range("T2").Value = 0.04784680267269592
$ws.Range("I3").Value = 0.0827560078572356
$ws.Range("J3").Value = 0.07616210790152332
$ws.Range("K3").Value = 0.08844667140905069
$ws.Range("L3").Value = 0.03829555421275314
$ws.Range("M3").Value = 0.076523636622613
$ws.Range("N3").Value = 0.07013559032659375
$ws.Range("O3").Value = 0.07829719358417983
$ws.Range("P3").Value = 0.035955920658406
$ws.Range("Q3").Value = 0.08314694127060565
$ws.Range("R3").Value = 0.06588454389756243
$ws.Range("S3").Value = 0.07346892209824178
$ws.Range("T3").Value = 0.04889327672913266
$ws.Range("I4").Value = 0.03430330156151285
$ws.Range("J4").Value = 0.03474366819625293
$ws.Range("K4").Value = 0.04356993503018398
$ws.Range("L4").Value = 0.04453060808782261
$ws.Range("M4").Value = 0.03691792969698193
$ws.Range("N4").Value = 0.0372511189117369
$ws.Range("O4").Value = 0.04725705504243226
$ws.Range("P4").Value = 0.04811176851273815
$ws.Range("Q4").Value = 0.03313946430072506
$ws.Range("R4").Value = 0.03373627663359698
$ws.Range("S4").Value = 0.04730205219343551
$ws.Range("T4").Value = 0.04726735691315499
$ws.Range("I5").Value = 0.09415936508721523
$ws.Range("J5").Value = 0.0963954712164433
$ws.Range("K5").Value = 0.1049889349555794
$ws.Range("L5").Value = 0.2160688388403541
$ws.Range("M5").Value = 0.1143676396755364
$ws.Range("N5").Value = 0.1153442466050348
$ws.Range("O5").Value = 0.1239450333778075
$ws.Range("P5").Value = 0.2084837940113204
$ws.Range("Q5").Value = 0.106304293606713
$ws.Range("R5").Value = 0.1082603464732572
$ws.Range("S5").Value = 0.1286414685093911
$ws.Range("T5").Value = 0.1725324567999956
$ws.Range("I6").Value = 0.04492404845963597
$ws.Range("J6").Value = 0.04614048166376553
$ws.Range("K6").Value = 0.05663656314976539
$ws.Range("L6").Value = 0.06076065697738244
$ws.Range("M6").Value = 0.05368638343389657
$ws.Range("N6").Value = 0.05433210740825254
$ws.Range("O6").Value = 0.07748465279634546
$ws.Range("P6").Value = 0.08146627715198362
$ws.Range("Q6").Value = 0.04875472538176546
$ws.Range("R6").Value = 0.04971623144676868
$ws.Range("S6").Value = 0.05940781452605793
$ws.Range("T6").Value = 0.06109338845926476
$ws.Range("I7").Value = 0.08088570093173948
$ws.Range("J7").Value = 0.08596300160112007
$ws.Range("K7").Value = 0.09299025208704255
$ws.Range("L7").Value = 0.2285944796822431
$ws.Range("M7").Value = 0.1009561020107799
$ws.Range("N7").Value = 0.1054672438089822
$ws.Range("O7").Value = 0.1203048575131326
$ws.Range("P7").Value = 0.2211003653838322
$ws.Range("Q7").Value = 0.08878711890512607
$ws.Range("R7").Value = 0.09246353879374736
$ws.Range("S7").Value = 0.09196459514394932
$ws.Range("T7").Value = 0.1318487156201701
$ws.Range("I8").Value = 0.7151896825054539
$ws.Range("J8").Value = 0.5982586464245429
$ws.Range("K8").Value = 0.8161034345445893
$ws.Range("L8").Value = -0.07323703676451936
$ws.Range("M8").Value = 0.6046697124422189
$ws.Range("N8").Value = 0.4913891232699088
$ws.Range("O8").Value = 0.6361205767952435
$ws.Range("P8").Value = -0.1147262582263156
$ws.Range("Q8").Value = 0.7221221882108834
$ws.Range("R8").Value = 0.4160044091078337
$ws.Range("S8").Value = 0.5504998074292227
$ws.Range("T8").Value = 0.11469463982465
$ws.Range("I9").Value = 0.2605515375973824
$ws.Range("J9").Value = 0.2156294237126781
$ws.Range("K9").Value = 0.3386792487804747
$ws.Range("L9").Value = -0.01509375344995455
$ws.Range("M9").Value = 0.1951876771888221
$ws.Range("N9").Value = 0.1586970759412157
$ws.Range("O9").Value = 0.242536423541914
$ws.Range("P9").Value = -0.02647535845312468
$ws.Range("Q9").Value = 0.2251154837217679
$ws.Range("R9").Value = 0.129636014327044
$ws.Range("S9").Value = 0.2024212792750599
$ws.Range("T9").Value = 0.03142198619997644
$ws.Range("I10").Value = 0.536469699503762
$ws.Range("J10").Value = 0.5341223780602828
$ws.Range("K10").Value = 0.4076293411010234
$ws.Range("L10").Value = 0.4199027069838933
$ws.Range("M10").Value = 0.519390806571452
$ws.Range("N10").Value = 0.5184992488839966
$ws.Range("O10").Value = 0.4995012972617772
$ws.Range("P10").Value = 0.5028079067912128
$ws.Range("Q10").Value = 0.5535658856193555
$ws.Range("R10").Value = 0.546271060586925
$ws.Range("S10").Value = 0.5069561373835163
$ws.Range("T10").Value = 0.5049264193526888
$ws.Range("I11").Value = 0.3774782903158234
$ws.Range("J11").Value = 0.3774782903158234
$ws.Range("K11").Value = 0.4020906920724777
$ws.Range("L11").Value = 0.4020906920724777
$ws.Range("M11").Value = 0.3805026795179351
$ws.Range("N11").Value = 0.3805026795179351
$ws.Range("O11").Value = 0.4039672598077267
$ws.Range("P11").Value = 0.4039672598077267
$ws.Range("Q11").Value = 0.3927447199451526
$ws.Range("R11").Value = 0.3927447199451526
$ws.Range("S11").Value = 0.4297621156472561
$ws.Range("T11").Value = 0.4297621156472561
$ws.Range("I12").Value = 0.2169835055494956
$ws.Range("J12").Value = 0.2169835055494956
$ws.Range("K12").Value = 0.2237973887343867
$ws.Range("L12").Value = 0.2237973887343867
$ws.Range("M12").Value = 0.2146205699109035
$ws.Range("N12").Value = 0.2146205699109035
$ws.Range("O12").Value = 0.2268436095699864
$ws.Range("P12").Value = 0.2268436095699864
$ws.Range("Q12").Value = 0.2415598577997207
$ws.Range("R12").Value = 0.2415598577997207
$ws.Range("S12").Value = 0.2263605968935193
$ws.Range("T12").Value = 0.2263605968935193
$ws.Range("I13").Value = 0.7324409187803218
$ws.Range("J13").Value = 0.7308367656736235
$ws.Range("K13").Value = 0.6384585664716415
$ws.Range("L13").Value = 0.647999002301619
$ws.Range("M13").Value = 0.7206877316643123
$ws.Range("N13").Value = 0.7200689195375651
$ws.Range("O13").Value = 0.706754057124384
$ws.Range("P13").Value = 0.7090894913839947
$ws.Range("Q13").Value = 0.7440200841505257
$ws.Range("R13").Value = 0.7391015225169852
$ws.Range("S13").Value = 0.7120085233924635
$ws.Range("T13").Value = 0.7105817471288497
$ws.Range("I14").Value = 35.37196486214975
$ws.Range("J14").Value = 35.37196486214975
$ws.Range("K14").Value = 368.4890547671815
$ws.Range("L14").Value = 368.4890547671815
$ws.Range("M14").Value = 38.00359170643648
$ws.Range("N14").Value = 38.00359170643648
$ws.Range("O14").Value = 338.8852096233285
$ws.Range("P14").Value = 338.8852096233285
$ws.Range("Q14").Value = 93.7022907473345
$ws.Range("R14").Value = 93.7022907473345
$ws.Range("S14").Value = 195.7782104719494
$ws.Range("T14").Value = 195.7782104719494

Write-Host "Edit complete"
